# cr1 completed and approved
# Append 7 more data rows (rows 9-15) to the "dataset" sheet, replicating the
# same repository record (tensorflow/ranking) with incrementing id values
# (8 through 14), matching the style/formatting of the last existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("dataset")

$lastRow = 8
$newRowsCount = 7

for ($i = 1; $i -le $newRowsCount; $i++) {
    $srcRow = $lastRow
    $dstRow = $lastRow + $i

    # Copy the whole source row (values + formatting) and paste into the new row
    $ws.Rows.Item($srcRow).Copy()
    $ws.Rows.Item($dstRow).PasteSpecial(-4104) | Out-Null

    # Update the incrementing id in column A (preserve the copied cell style)
    $ws.Cells.Item($dstRow, 1).Value = 7 + $i
    $ws.Cells.Item($dstRow, 1).Style = $ws.Cells.Item($srcRow, 1).Style
}

$excel.CutCopyMode = 0

# Refresh the sheet dimension to cover the newly added rows
$ws.UsedRange | Out-Null
